$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---
# Row 3 corresponds to af843c2a-ac76-495f-a669-548f82605fb1.md which has now
# been handed back (in sync with en-US) for both locales.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet ---
# Row 3 (af843c2a...) status flips to handed-back, and gets a fresh
# handback timestamp.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("H3").Value = "2016-03-19 22:37:18"

# --- de-de sheet ---
# Same update for the de-de locale, with its own handback timestamp.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("H3").Value = "2016-03-19 22:37:23"
